$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: "001" -> "002" (must remain text, not be auto-converted to a number)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("J2").ClearFormats()

# N2: report date text changes
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Numeric financial figures for row 2
$ws.Range("O2").Value = 123530774.88
$ws.Range("P2").Value = 494606220.74
$ws.Range("Q2").Value = 344920639.83
$ws.Range("S2").Value = 277012678.26
$ws.Range("T2").Value = 277012678.26
$ws.Range("V2").Value = 16430139.25
$ws.Range("W2").Value = 21485922.61
$ws.Range("X2").Value = -1564916.87
$ws.Range("Y2").Value = 141941697.98
$ws.Range("Z2").Value = 141510189.61
$ws.Range("AA2").Value = 17979414.73
$ws.Range("AG2").Value = 1694018.85
$ws.Range("AS2").Value = 122586674.88
